$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("I2").Value = 0.4458169960283037
$ws.Range("J2").Value = 0.4458169960283037
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 84.16319028771844
$ws.Range("R2").Value = 757.468712589466
$ws.Range("S2").Value = 0.1705085010928244
$ws.Range("T2").Value = 0.1705085010928244
$ws.Range("I3").Value = 0.4458169960283037
$ws.Range("J3").Value = 0.4458169960283037
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("S3").Value = 0.1455326065497182
$ws.Range("T3").Value = 0.1455326065497182
$ws.Range("I4").Value = 0.4458169960283037
$ws.Range("J4").Value = 0.4458169960283037
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 21.40102087590244
$ws.Range("R4").Value = 192.609187883122
$ws.Range("S4").Value = 0.04335691148270149
$ws.Range("T4").Value = 0.04335691148270147
$ws.Range("I5").Value = 0.4458169960283037
$ws.Range("J5").Value = 0.4458169960283037
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 42.65650539970778
$ws.Range("R5").Value = 383.90854859737
$ws.Range("S5").Value = 0.08641897690305954
$ws.Range("T5").Value = 0.08641897690305951
$ws.Range("G6").Value = 1.399743666666667
$ws.Range("H6").Value = 4.199231
$ws.Range("I6").Value = 0.2598558798146963
$ws.Range("J6").Value = 0.2598558798146962
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 49.05667584471878
$ws.Range("R6").Value = 441.5100826024691
$ws.Range("S6").Value = 0.09938525664586378
$ws.Range("T6").Value = 0.09938525664586374
$ws.Range("G7").Value = 1.399743666666667
$ws.Range("H7").Value = 4.199231
$ws.Range("I7").Value = 0.2598558798146963
$ws.Range("J7").Value = 0.2598558798146962
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("Q7").Value = 41.87090883204633
$ws.Range("R7").Value = 376.838179488417
$ws.Range("S7").Value = 0.08482741540500204
$ws.Range("T7").Value = 0.084827415405002
$ws.Range("G8").Value = 1.399743666666667
$ws.Range("H8").Value = 4.199231
$ws.Range("I8").Value = 0.2598558798146963
$ws.Range("J8").Value = 0.2598558798146962
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 12.47413436047478
$ws.Range("R8").Value = 112.267209244273
$ws.Range("S8").Value = 0.02527168878655766
$ws.Range("T8").Value = 0.02527168878655765
$ws.Range("G9").Value = 1.399743666666667
$ws.Range("H9").Value = 4.199231
$ws.Range("I9").Value = 0.2598558798146963
$ws.Range("J9").Value = 0.2598558798146962
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 24.86343912235611
$ws.Range("R9").Value = 223.770952101205
$ws.Range("S9").Value = 0.05037151897727279
$ws.Range("T9").Value = 0.05037151897727278
$ws.Range("G10").Value = 1.585427
$ws.Range("H10").Value = 4.756281
$ws.Range("I10").Value = 0.294327124157
$ws.Range("J10").Value = 0.294327124157
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 55.56430099782433
$ws.Range("R10").Value = 500.078708980419
$ws.Range("S10").Value = 0.1125692318105019
$ws.Range("T10").Value = 0.1125692318105018
$ws.Range("G11").Value = 1.585427
$ws.Range("H11").Value = 4.756281
$ws.Range("I11").Value = 0.294327124157
$ws.Range("J11").Value = 0.294327124157
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 47.425304330863
$ws.Range("R11").Value = 426.827738977767
$ws.Range("S11").Value = 0.0960802166325021
$ws.Range("T11").Value = 0.09608021663250207
$ws.Range("G12").Value = 1.585427
$ws.Range("H12").Value = 4.756281
$ws.Range("I12").Value = 0.294327124157
$ws.Range("J12").Value = 0.294327124157
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 14.12889365938033
$ws.Range("R12").Value = 127.160042934423
$ws.Range("S12").Value = 0.02862411075109162
$ws.Range("T12").Value = 0.02862411075109161
$ws.Range("G13").Value = 1.585427
$ws.Range("H13").Value = 4.756281
$ws.Range("I13").Value = 0.294327124157
$ws.Range("J13").Value = 0.294327124157
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 28.16170462932833
$ws.Range("R13").Value = 253.455341663955
$ws.Range("S13").Value = 0.05705356496290441
$ws.Range("T13").Value = 0.05705356496290439
